# 27.12.2021 refactor class EI and Chrome
# The BOM sheet's "ParentId" column (A) for rows 2-6 is updated to a new
# part code. Three new codes were appended to the shared-string pool by the
# original tool; the last of them ("B12UB51111212") is the one actually
# used to replace the previous value ("B12UB21111212") in A2:A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParentId = "B12UB51111212"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newParentId
}
